# Auto-generated Excel COM-interop edit script
# Applies updated market-price-derived profit figures to the Leve profit
# tracking sheets (ALC, ARM, BSM, CUL, GSM, LTW, WVR) as produced by the
# scheduled price-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 499.8
$ws.Range("I6").Value = 533.6667
$ws.Range("J6").Value = 195
$ws.Range("K6").Value = 1601.0001
$ws.Range("L6").Value = 585
$ws.Range("M6").Value = -1489.0001
$ws.Range("N6").Value = -809

$ws.Range("H33").Value = 409.75
$ws.Range("I33").Value = 421.7
$ws.Range("J33").Value = 350
$ws.Range("K33").Value = 421.7
$ws.Range("L33").Value = 350
$ws.Range("M33").Value = -192.7
$ws.Range("N33").Value = -808

$ws.Range("H62").Value = 4167.2666
$ws.Range("J62").Value = 5143.7144
$ws.Range("L62").Value = 5143.7144
$ws.Range("N62").Value = -6391.7144

$ws.Range("H65").Value = 4167.2666
$ws.Range("J65").Value = 5143.7144
$ws.Range("L65").Value = 25718.572
$ws.Range("N65").Value = -31958.572

$ws.Range("H74").Value = 37833.875
$ws.Range("I74").Value = 39044.043
$ws.Range("K74").Value = 39044.043
$ws.Range("M74").Value = -38108.043

$ws.Range("H76").Value = 6678.3125
$ws.Range("I76").Value = 5850.909
$ws.Range("K76").Value = 5850.909
$ws.Range("M76").Value = -5535.909

$ws.Range("H77").Value = 37833.875
$ws.Range("I77").Value = 39044.043
$ws.Range("K77").Value = 195220.215
$ws.Range("M77").Value = -190540.215

$ws.Range("H79").Value = 6678.3125
$ws.Range("I79").Value = 5850.909
$ws.Range("K79").Value = 5850.909
$ws.Range("M79").Value = -4758.909

$ws.Range("H106").Value = 3198.6316
$ws.Range("I106").Value = 2618.4666
$ws.Range("K106").Value = 2618.4666
$ws.Range("M106").Value = -1987.4666

$ws.Range("H137").Value = 1753.9048
$ws.Range("I137").Value = 1257.375
$ws.Range("J137").Value = 3342.8
$ws.Range("K137").Value = 3772.125
$ws.Range("L137").Value = 10028.4
$ws.Range("M137").Value = -1222.125
$ws.Range("N137").Value = -15128.4

$ws.Range("H138").Value = 3100.204
$ws.Range("I138").Value = 2689.85
$ws.Range("J138").Value = 3383.2068
$ws.Range("K138").Value = 8069.549999999999
$ws.Range("L138").Value = 10149.6204
$ws.Range("M138").Value = -2929.549999999999
$ws.Range("N138").Value = -20429.6204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 11927.25
$ws.Range("J24").Value = 11927.25
$ws.Range("L24").Value = 11927.25
$ws.Range("N24").Value = -12675.25

$ws.Range("H32").Value = 3833.926
$ws.Range("I32").Value = 2780.9
$ws.Range("K32").Value = 2780.9
$ws.Range("M32").Value = -2493.9

$ws.Range("H61").Value = 37375250
$ws.Range("I61").Value = 37375250
$ws.Range("K61").Value = 37375250
$ws.Range("M61").Value = -37375038

$ws.Range("H88").Value = 2761.5
$ws.Range("I88").Value = 2300
$ws.Range("J88").Value = 2915.3333
$ws.Range("K88").Value = 2300
$ws.Range("L88").Value = 2915.3333
$ws.Range("M88").Value = -1894
$ws.Range("N88").Value = -3727.3333

$ws.Range("H91").Value = 2761.5
$ws.Range("I91").Value = 2300
$ws.Range("J91").Value = 2915.3333
$ws.Range("K91").Value = 2300
$ws.Range("L91").Value = 2915.3333
$ws.Range("M91").Value = -896
$ws.Range("N91").Value = -5723.3333

$ws.Range("H100").Value = 11927.25
$ws.Range("J100").Value = 11927.25
$ws.Range("L100").Value = 11927.25
$ws.Range("N100").Value = -14091.25

$ws.Range("H110").Value = 79057.92
$ws.Range("I110").Value = 92249.27
$ws.Range("J110").Value = 6505.5
$ws.Range("K110").Value = 92249.27
$ws.Range("L110").Value = 6505.5
$ws.Range("M110").Value = -90204.27
$ws.Range("N110").Value = -10595.5

$ws.Range("H136").Value = 37375250
$ws.Range("I136").Value = 37375250
$ws.Range("K136").Value = 112125750
$ws.Range("M136").Value = -112123200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4080.5
$ws.Range("I86").Value = 2050
$ws.Range("K86").Value = 2050
$ws.Range("M86").Value = -927

$ws.Range("H89").Value = 4080.5
$ws.Range("I89").Value = 2050
$ws.Range("K89").Value = 10250
$ws.Range("M89").Value = -4634

$ws.Range("H99").Value = 1388.2727
$ws.Range("I99").Value = 826
$ws.Range("K99").Value = 826
$ws.Range("M99").Value = 672

$ws.Range("H105").Value = 1920.2858
$ws.Range("I105").Value = 1407.1666
$ws.Range("K105").Value = 1407.1666
$ws.Range("M105").Value = 339.8334

$ws.Range("H132").Value = 134999.5
$ws.Range("J132").Value = 134999.5
$ws.Range("L132").Value = 134999.5
$ws.Range("N132").Value = -145119.5

$ws.Range("H134").Value = 10640830
$ws.Range("I134").Value = 11906644
$ws.Range("K134").Value = 35719932
$ws.Range("M134").Value = -35717397

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2779.8
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 3224.75
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 9674.25
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -11670.25

$ws.Range("H78").Value = 2779.8
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 3224.75
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 29022.75
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -39006.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10003.523
$ws.Range("I70").Value = 9275.714
$ws.Range("K70").Value = 9275.714
$ws.Range("M70").Value = -9005.714

$ws.Range("H73").Value = 10003.523
$ws.Range("I73").Value = 9275.714
$ws.Range("K73").Value = 9275.714
$ws.Range("M73").Value = -8339.714

$ws.Range("H80").Value = 3532.3333
$ws.Range("I80").Value = 3749.5
$ws.Range("K80").Value = 3749.5
$ws.Range("M80").Value = -2751.5

$ws.Range("H83").Value = 3532.3333
$ws.Range("I83").Value = 3749.5
$ws.Range("K83").Value = 18747.5
$ws.Range("M83").Value = -13755.5

$ws.Range("H113").Value = 56339
$ws.Range("I113").Value = 93598.09
$ws.Range("K113").Value = 93598.09
$ws.Range("M113").Value = -91428.09

$ws.Range("H122").Value = 83359.266
$ws.Range("I122").Value = 87527.78999999999
$ws.Range("K122").Value = 262583.37
$ws.Range("M122").Value = -260133.37

$ws.Range("H132").Value = 7813664.5
$ws.Range("I132").Value = 8334442
$ws.Range("K132").Value = 25003326
$ws.Range("M132").Value = -25000796

$ws.Range("H138").Value = 111796.336
$ws.Range("I138").Value = 110390
$ws.Range("K138").Value = 110390
$ws.Range("M138").Value = -105250

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents() | Out-Null

$ws.Range("H132").Value = 48013210
$ws.Range("I132").Value = 68588620
$ws.Range("J132").Value = 3899.3333
$ws.Range("K132").Value = 205765860
$ws.Range("L132").Value = 11697.9999
$ws.Range("M132").Value = -205763330
$ws.Range("N132").Value = -16757.9999

$ws.Range("H136").Value = 2465
$ws.Range("I136").Value = 2197.75
$ws.Range("K136").Value = 6593.25
$ws.Range("M136").Value = -4043.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1875.56
$ws.Range("I96").Value = 1681.091
$ws.Range("J96").Value = 2028.3572
$ws.Range("K96").Value = 1681.091
$ws.Range("L96").Value = 2028.3572
$ws.Range("M96").Value = -308.0909999999999
$ws.Range("N96").Value = -4774.3572
